$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data fix: rows 98-127 -------------------------------------------------
# Column F (Drug_Dose) was 50 for every one of these rows and should be 100.
# Column I (Dose_Frequency) was "BID" and should become "QD".
for ($r = 98; $r -le 127; $r++) {
    $ws.Cells.Item($r, 6).Value = 100
    $ws.Cells.Item($r, 9).Value = "QD"
}

# --- View state: scroll position & active selection ------------------------
# The sheet view had scrolled to A219 with D246 selected; move it to show
# B80 in the top-left corner with I93 selected instead.
$excel.ActiveWindow.ScrollRow = 80
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("I93").Select()
